# Apply "Automatic update of files" changes to the EKERÖ overview sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") on every data row (2-45) moves forward one day:
# 45189 (2023-09-20) -> 45190 (2023-09-21)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}

# Row 2 additionally gets an extra observed species ("Spricktaggsvamp"),
# which bumps the VU / Rödlistade / Hotade / Alla arter counters.
$ws.Range("K2").Value = 1
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = "Ryl`r`nSpricktaggsvamp`r`nBacktimjan`r`nOrange taggsvamp`r`nSpillkråka`r`nDropptaggsvamp`r`nGrönpyrola`r`nKopparödla"
